$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (shifting old Quantity/Price/Grant Date columns
# J,K,L one place to the right, to K,L,M) so there is room for the new
# "Preferred Conversion" column.
$xlShiftToRight = -4161
$ws.Range("J1:J7").Insert($xlShiftToRight)

# New header for the inserted column.
$ws.Range("J1").Value = "Preferred Conversion"

# Populate per-row "Preferred Conversion" values (blank where not applicable).
$ws.Range("J3").Value = 2
$ws.Range("J5").Value = 3

# Update the active cell selection to match the edited sheet's state.
$ws.Range("J2").Select()
